$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "skiprow"
$ws.Range("B1").Value = "template"
$ws.Range("C1").Value = "to"
$ws.Range("D1").Value = "from"
$ws.Range("E1").Value = "redirect"
$ws.Range("F1").Value = "name"
$ws.Range("F1").Select()
